$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data for "impactParameter" / "transit impact parameter"
$ws.Range("A11").Value = "impactParameter"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = "m"
$ws.Range("D11").Value = "transit impact parameter"
$ws.Range("E11").Value = "to be modified based on the data being used, so maybe it shouldn't be in here"

# Update the selected cell to match the post-edit state
$ws.Range("B12").Select()
